$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy formatting (number format/style) from the new column's right neighbor (E)
# into the newly inserted (blank) column D, for all relevant rows
$src = $ws.Range("E5:E102")
$dst = $ws.Range("D5:D102")
$src.Copy()
$dst.PasteSpecial(-4122)

# Rows 5, 6, 37 and 79 only have a single label cell (in column A or B) and never
# had any cell in column D/E - remove the blank styled cell that PasteSpecial created
# so the structure matches rows that truly have no data in that column.
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# Populate the newly inserted column D with the latest quarter figures,
# and update the handful of prior-quarter cells that were restated.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 487300
$ws.Range("G8").Value = 429900
$ws.Range("H8").Value = 472900
$ws.Range("D9").Value = 258000
$ws.Range("G9").Value = 244200
$ws.Range("H9").Value = 258100
$ws.Range("D10").Value = 229300
$ws.Range("G10").Value = 185700
$ws.Range("H10").Value = 214800
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 2900
$ws.Range("D15").Value = 35800
$ws.Range("G15").Value = 72500
$ws.Range("D17").Value = 395500
$ws.Range("G17").Value = 401200
$ws.Range("H17").Value = 401400
$ws.Range("D18").Value = 91800
$ws.Range("G18").Value = 28700
$ws.Range("H18").Value = 71500
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 127500
$ws.Range("G21").Value = 64900
$ws.Range("H21").Value = 100600
$ws.Range("D22").Value = 34100
$ws.Range("D23").Value = 57700
$ws.Range("G23").Value = -1700
$ws.Range("H23").Value = 42300
$ws.Range("D24").Value = 20000
$ws.Range("G24").Value = -90700
$ws.Range("H24").Value = 12600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 37700
$ws.Range("G26").Value = 88900
$ws.Range("H26").Value = 29700
$ws.Range("D27").Value = 37700
$ws.Range("G27").Value = 88900
$ws.Range("H27").Value = 29700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("F29").Value = 110000
$ws.Range("G29").Value = 1400
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 37700
$ws.Range("G33").Value = 198900
$ws.Range("H33").Value = 31100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 37700
$ws.Range("G35").Value = 198900
$ws.Range("H35").Value = 31100
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 164100
$ws.Range("D42").Value = 17100
$ws.Range("D43").Value = 42200
$ws.Range("D44").Value = 15300
$ws.Range("D45").Value = 187600
$ws.Range("D46").Value = 426200
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 1488900
$ws.Range("D49").Value = 1909700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 652400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4477200
$ws.Range("D57").Value = 33300
$ws.Range("D58").Value = 400
$ws.Range("D59").Value = 228200
$ws.Range("D60").Value = 261900
$ws.Range("D61").Value = 2967400
$ws.Range("D62").Value = 211700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3441000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 290300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1036200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 37700
$ws.Range("G81").Value = 198900
$ws.Range("H81").Value = 31100
$ws.Range("D83").Value = 35800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 114900
$ws.Range("D91").Value = -33900
$ws.Range("G91").Value = -30600
$ws.Range("H91").Value = -23100
$ws.Range("I91").Value = -23600
$ws.Range("J91").Value = -6200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -325800
$ws.Range("G94").Value = -41800
$ws.Range("H94").Value = -23100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 774000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 563100
$ws.Range("G102").Value = -29100
$ws.Range("H102").Value = 30000

Write-Host "Edit complete"
